# Update statistics on the "2o Parcial" and "Final" sheets
# (Estadisticos Segundo Parcial 26 Mayo)

$wb = $excel.ActiveWorkbook

# --- Sheet "2o Parcial" ---
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Row 8
$ws2.Range("E8").Value = 14
$ws2.Range("F8").Value = 1
$ws2.Range("G8").Value = 93.3
$ws2.Range("H8").Value = 6.7
$ws2.Range("I8").Value = 8.300000000000001
$ws2.Range("J8").Value = 0
$ws2.Range("K8").Value = 0

# Row 9
$ws2.Range("E9").Value = 24
$ws2.Range("F9").Value = 1
$ws2.Range("G9").Value = 96
$ws2.Range("H9").Value = 4
$ws2.Range("I9").Value = 8.4
$ws2.Range("J9").Value = 0
$ws2.Range("K9").Value = 0

# Row 10
$ws2.Range("E10").Value = 38
$ws2.Range("F10").Value = 2
$ws2.Range("G10").Value = 95
$ws2.Range("H10").Value = 5
$ws2.Range("I10").Value = 8.4
$ws2.Range("J10").Value = 0
$ws2.Range("K10").Value = 0

# Row 20 (Totales Generales)
$ws2.Range("E20").Value = 360
$ws2.Range("F20").Value = 36
$ws2.Range("G20").Value = 90.90000000000001
$ws2.Range("H20").Value = 9.1
$ws2.Range("I20").Value = 8.300000000000001
$ws2.Range("J20").Value = 0
$ws2.Range("K20").Value = 0

# --- Sheet "Final" ---
$ws3 = $wb.Worksheets.Item("Final")

# Row 8
$ws3.Range("I8").Value = 8.199999999999999

# Row 9
$ws3.Range("E9").Value = 24
$ws3.Range("F9").Value = 1
$ws3.Range("G9").Value = 96
$ws3.Range("H9").Value = 4
$ws3.Range("I9").Value = 8.4

# Row 10
$ws3.Range("E10").Value = 38
$ws3.Range("F10").Value = 2
$ws3.Range("G10").Value = 95
$ws3.Range("H10").Value = 5
$ws3.Range("I10").Value = 8.300000000000001

# Row 20 (Totales Generales)
$ws3.Range("E20").Value = 360
$ws3.Range("F20").Value = 36
$ws3.Range("G20").Value = 90.90000000000001
$ws3.Range("H20").Value = 9.1
$ws3.Range("I20").Value = 8.300000000000001
